# Generate Report for Handback
# The row for file "91d9f6e7-7253-4b94-8469-5d80cf80535c" has been handed back
# successfully (in sync with en-US). Update the Overview sheet and the two
# per-locale sheets (zh-cn, de-de) to reflect the new status and a refreshed
# "Latest Handback DateTime" for the zh-cn locale, and clear the previous
# error detail message.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$newHandbackDateTime = "2016-11-09 00:43:21"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value = $newStatus
$overview.Range("F6").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = $newStatus
$zhcn.Range("K6").Value = $newHandbackDateTime
$zhcn.Range("P6").Value = ""

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = $newStatus
$dede.Range("P6").Value = ""
